# Auto-generated edit script applying the diff's cell value changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2313703
$ws.Range("J17").Value = 2313703
$ws.Range("L17").Value = 6941109
$ws.Range("N17").Value = -6941445
$ws.Range("H38").Value = 920.6667
$ws.Range("I38").Value = 150.5
$ws.Range("J38").Value = 1140.7142
$ws.Range("K38").Value = 451.5
$ws.Range("L38").Value = 3422.1426
$ws.Range("M38").Value = -79.5
$ws.Range("N38").Value = -4166.142599999999
$ws.Range("H94").Value = 3666.3333
$ws.Range("I94").Value = 3666.3333
$ws.Range("K94").Value = 3666.3333
$ws.Range("M94").Value = -3215.3333
$ws.Range("H107").Value = 39061.77
$ws.Range("I107").Value = 45891.184
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 45891.184
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -43971.184
$ws.Range("N107").Value = -5340
$ws.Range("H111").Value = 806.2857
$ws.Range("I111").Value = 564.5333000000001
$ws.Range("J111").Value = 1410.6666
$ws.Range("K111").Value = 1693.5999
$ws.Range("L111").Value = 4231.9998
$ws.Range("M111").Value = 1373.4001
$ws.Range("N111").Value = -10365.9998
$ws.Range("H132").Value = 7869.587
$ws.Range("I132").Value = 5912.8604
$ws.Range("J132").Value = 35916
$ws.Range("K132").Value = 17738.5812
$ws.Range("L132").Value = 107748
$ws.Range("M132").Value = -15208.5812
$ws.Range("N132").Value = -112808
$ws.Range("H138").Value = 4015.5952
$ws.Range("I138").Value = 4269.6523
$ws.Range("K138").Value = 12808.9569
$ws.Range("M138").Value = -7668.956899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1856.6666
$ws.Range("I2").Value = 635
$ws.Range("K2").Value = 635
$ws.Range("M2").Value = -522
$ws.Range("H5").Value = 180.9
$ws.Range("I5").Value = 180.9
$ws.Range("K5").Value = 180.9
$ws.Range("M5").Value = -68.90000000000001
$ws.Range("H32").Value = 12556.696
$ws.Range("I32").Value = 4237.7344
$ws.Range("K32").Value = 4237.7344
$ws.Range("M32").Value = -3950.7344
$ws.Range("H37").Value = 8756
$ws.Range("I37").Value = 8756
$ws.Range("K37").Value = 8756
$ws.Range("M37").Value = -8483
$ws.Range("H97").Value = 2585.6667
$ws.Range("I97").Value = 906.7826
$ws.Range("J97").Value = 8102
$ws.Range("K97").Value = 906.7826
$ws.Range("L97").Value = 8102
$ws.Range("M97").Value = -410.7826
$ws.Range("N97").Value = -9094
$ws.Range("H116").Value = 1856.6666
$ws.Range("I116").Value = 635
$ws.Range("K116").Value = 635
$ws.Range("M116").Value = 1659
$ws.Range("H132").Value = 1593633.1
$ws.Range("I132").Value = 2579.0952
$ws.Range("K132").Value = 7737.285600000001
$ws.Range("M132").Value = -5207.285600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1856.6666
$ws.Range("I3").Value = 635
$ws.Range("K3").Value = 635
$ws.Range("M3").Value = -521
$ws.Range("H4").Value = 180.9
$ws.Range("I4").Value = 180.9
$ws.Range("K4").Value = 180.9
$ws.Range("M4").Value = -65.90000000000001
$ws.Range("H22").Value = 540.375
$ws.Range("I22").Value = 540.375
$ws.Range("K22").Value = 540.375
$ws.Range("M22").Value = -367.375
$ws.Range("H134").Value = 11233.471
$ws.Range("I134").Value = 5076.44
$ws.Range("K134").Value = 15229.32
$ws.Range("M134").Value = -12694.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 200603
$ws.Range("I10").Value = 754
$ws.Range("J10").Value = 999999
$ws.Range("K10").Value = 754
$ws.Range("L10").Value = 999999
$ws.Range("M10").Value = -615
$ws.Range("N10").Value = -1000277
$ws.Range("H16").Value = 3530.0476
$ws.Range("I16").Value = 1445
$ws.Range("J16").Value = 6918.25
$ws.Range("K16").Value = 1445
$ws.Range("L16").Value = 6918.25
$ws.Range("M16").Value = -1158
$ws.Range("N16").Value = -7492.25
$ws.Range("H31").Value = 17908.77
$ws.Range("J31").Value = 21145.371
$ws.Range("L31").Value = 21145.371
$ws.Range("N31").Value = -21735.371
$ws.Range("H34").Value = 17908.77
$ws.Range("J34").Value = 21145.371
$ws.Range("L34").Value = 21145.371
$ws.Range("N34").Value = -21549.371
$ws.Range("H58").Value = 13175.149
$ws.Range("I58").Value = 5883.294
$ws.Range("K58").Value = 5883.294
$ws.Range("M58").Value = -5680.294
$ws.Range("H105").Value = 18094.889
$ws.Range("J105").Value = 10216.5
$ws.Range("L105").Value = 10216.5
$ws.Range("N105").Value = -13710.5
$ws.Range("H113").Value = 3530.0476
$ws.Range("I113").Value = 1445
$ws.Range("J113").Value = 6918.25
$ws.Range("K113").Value = 1445
$ws.Range("L113").Value = 6918.25
$ws.Range("M113").Value = 725
$ws.Range("N113").Value = -11258.25
$ws.Range("H122").Value = 6777.636
$ws.Range("I122").Value = 4414.3125
$ws.Range("K122").Value = 13242.9375
$ws.Range("M122").Value = -10792.9375
$ws.Range("H132").Value = 7322.7334
$ws.Range("I132").Value = 3585.2
$ws.Range("J132").Value = 9191.5
$ws.Range("K132").Value = 10755.6
$ws.Range("L132").Value = 27574.5
$ws.Range("M132").Value = -8225.599999999999
$ws.Range("N132").Value = -32634.5
$ws.Range("H136").Value = 13175.149
$ws.Range("I136").Value = 5883.294
$ws.Range("K136").Value = 17649.882
$ws.Range("M136").Value = -15099.882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 366.77777
$ws.Range("J2").Value = 303.66666
$ws.Range("L2").Value = 1821.99996
$ws.Range("N2").Value = -2047.99996
$ws.Range("H38").Value = 75.44444
$ws.Range("I38").Value = 19.8
$ws.Range("J38").Value = 145
$ws.Range("K38").Value = 59.40000000000001
$ws.Range("L38").Value = 435
$ws.Range("M38").Value = 287.6
$ws.Range("N38").Value = -1129
$ws.Range("H131").Value = 1498.4
$ws.Range("J131").Value = 1499.9375
$ws.Range("L131").Value = 4499.8125
$ws.Range("N131").Value = -14579.8125
$ws.Range("H132").Value = 1538.5555
$ws.Range("I132").Value = 1246.3334
$ws.Range("J132").Value = 1684.6666
$ws.Range("K132").Value = 11217.0006
$ws.Range("L132").Value = 15161.9994
$ws.Range("M132").Value = -8687.000599999999
$ws.Range("N132").Value = -20221.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5319.5454
$ws.Range("I122").Value = 3641.35
$ws.Range("J122").Value = 7901.385
$ws.Range("K122").Value = 10924.05
$ws.Range("L122").Value = 23704.155
$ws.Range("M122").Value = -8474.049999999999
$ws.Range("N122").Value = -28604.155
$ws.Range("H132").Value = 5557.147
$ws.Range("I132").Value = 1866.4546
$ws.Range("K132").Value = 5599.3638
$ws.Range("M132").Value = -3069.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17483.084
$ws.Range("I22").Value = 15310.777
$ws.Range("K22").Value = 15310.777
$ws.Range("M22").Value = -15015.777
$ws.Range("H27").Value = 17483.084
$ws.Range("I27").Value = 15310.777
$ws.Range("K27").Value = 15310.777
$ws.Range("M27").Value = -15203.777
$ws.Range("H39").Value = 5830
$ws.Range("J39").Value = 4995
$ws.Range("L39").Value = 4995
$ws.Range("N39").Value = -5915
$ws.Range("H40").Value = 9794.789000000001
$ws.Range("I40").Value = 7206.8
$ws.Range("K40").Value = 7206.8
$ws.Range("M40").Value = -7070.8
$ws.Range("H93").Value = 13409.471
$ws.Range("I93").Value = 10198.5
$ws.Range("J93").Value = 16263.667
$ws.Range("K93").Value = 10198.5
$ws.Range("L93").Value = 16263.667
$ws.Range("M93").Value = -8950.5
$ws.Range("N93").Value = -18759.667
$ws.Range("H122").Value = 6456.1353
$ws.Range("I122").Value = 4537.88
$ws.Range("K122").Value = 13613.64
$ws.Range("M122").Value = -11163.64

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 25383.75
$ws.Range("I51").Value = 19356.666
$ws.Range("J51").Value = 29000
$ws.Range("K51").Value = 19356.666
$ws.Range("L51").Value = 29000
$ws.Range("M51").Value = -18846.666
$ws.Range("N51").Value = -30020
$ws.Range("H54").Value = 31035
$ws.Range("I54").Value = 31035
$ws.Range("K54").Value = 31035
$ws.Range("M54").Value = -30515
$ws.Range("H64").Value = 52954.117
$ws.Range("J64").Value = 52957.5
$ws.Range("L64").Value = 52957.5
$ws.Range("N64").Value = -53453.5
$ws.Range("H67").Value = 52954.117
$ws.Range("J67").Value = 52957.5
$ws.Range("L67").Value = 52957.5
$ws.Range("N67").Value = -54673.5
$ws.Range("H75").Value = 34989
$ws.Range("I75").Value = 34989
$ws.Range("K75").Value = 34989
$ws.Range("M75").Value = -34053
$ws.Range("H78").Value = 34989
$ws.Range("I78").Value = 34989
$ws.Range("K78").Value = 104967
$ws.Range("M78").Value = -100287
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 4859.516
$ws.Range("I122").Value = 1782.6
$ws.Range("K122").Value = 5347.799999999999
$ws.Range("M122").Value = -2897.799999999999
$ws.Range("H132").Value = 5136.4194
$ws.Range("I132").Value = 1777.5366
$ws.Range("K132").Value = 5332.6098
$ws.Range("M132").Value = -2802.6098
